# Regenerate save_data to use K (Strike count) instead of Strike# column,
# updating the K column values based on recalculated std/mean s_vals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new K (column G) value
$updates = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 1
    7  = 2
    9  = 1
    11 = 2
    12 = 1
    14 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
